$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 1.53
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.1
$ws.Range("K4").Value = 2.38
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 15
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("U4").Value = 1.33
$ws.Range("V4").Value = 3.25
$ws.Range("W4").Value = 1.73
$ws.Range("X4").Value = 2
$ws.Range("Z4").Value = 8
$ws.Range("AB4").Value = 12
$ws.Range("AF4").Value = 8
$ws.Range("AG4").Value = 15
$ws.Range("AI4").Value = 17
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 51

# Row 6 updates
$ws.Range("M6").Value = 1.03
$ws.Range("O6").Value = 1.19
$ws.Range("T6").Value = 1.37
